# Applies a weekly reshuffle of the "Fecha" (date) and associated price/volume
# columns across the data rows of the sheet. Each destination row's D, M, N, O,
# P, R and S values are replaced by the values that currently live in a
# corresponding source row (a permutation of the existing rows). Row 10 is not
# part of any cycle and stays as-is.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# destination row -> source row (values are copied from source row's
# original content into the destination row)
$mapping = @{
    2  = 13
    3  = 15
    4  = 5
    5  = 11
    6  = 16
    7  = 9
    8  = 7
    9  = 8
    11 = 4
    12 = 17
    13 = 3
    14 = 6
    15 = 2
    16 = 12
    17 = 14
}

$columns = @("D", "M", "N", "O", "P", "R", "S")

# Snapshot the current ("before") values of every column we touch, for every
# row involved, before any writes happen. Use Value2 (not Value) to get the
# raw underlying value rather than a wrapped Variant object.
$snapshot = @{}
$allRows = @()
foreach ($row in $mapping.Keys) { $allRows += $row }
foreach ($row in $mapping.Values) { $allRows += $row }
$allRows = $allRows | Select-Object -Unique

foreach ($row in $allRows) {
    $rowValues = @{}
    foreach ($col in $columns) {
        $rowValues[$col] = $ws.Range("$col$row").Value2
    }
    $snapshot[$row] = $rowValues
}

# Now write the destination rows using the snapshotted source-row values.
foreach ($destRow in $mapping.Keys) {
    $srcRow = $mapping[$destRow]
    $srcValues = $snapshot[$srcRow]
    foreach ($col in $columns) {
        $ws.Range("$col$destRow").Value2 = $srcValues[$col]
    }
}
